$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B12 value (hours worked entry changed from 0 to 9)
$ws.Range("B12").Value = 9

# Update the active selection to B12 to match the saved cursor position
$ws.Range("B12").Select()
